$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (60) down onto the two
# new rows (61, 62) so the new cells pick up the same styles (bold/border
# index for column A, date-time number format for column E) without
# introducing any new style entries.
$ws.Range("A60:V60").Copy()
$ws.Range("A61:V62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Row 61: Sevilla 5 - 1 Almeria ----
$ws.Range("A61").Value = 60
$ws.Range("B61").Value = "spain"
$ws.Range("C61").Value = "laliga"
$ws.Range("D61").Value = "2023-2024"
$ws.Range("E61").Value = 45195.79166666666
$ws.Range("F61").Value = "Sevilla"
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = "Almeria"
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 1.67
$ws.Range("K61").Value = "16/09/2023 21:01"
$ws.Range("L61").Value = 1.71
$ws.Range("M61").Value = "26/09/2023 18:57"
$ws.Range("N61").Value = 4.08
$ws.Range("O61").Value = "16/09/2023 21:01"
$ws.Range("P61").Value = 3.97
$ws.Range("Q61").Value = "26/09/2023 18:56"
$ws.Range("R61").Value = 5.18
$ws.Range("S61").Value = "16/09/2023 21:01"
$ws.Range("T61").Value = 5.18
$ws.Range("U61").Value = "26/09/2023 18:57"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/spain/laliga/sevilla-almeria/EJCPWy4o/"

# ---- Row 62: Mallorca 2 - 2 Barcelona ----
$ws.Range("A62").Value = 61
$ws.Range("B62").Value = "spain"
$ws.Range("C62").Value = "laliga"
$ws.Range("D62").Value = "2023-2024"
$ws.Range("E62").Value = 45195.89583333334
$ws.Range("F62").Value = "Mallorca"
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = "Barcelona"
$ws.Range("I62").Value = 2
$ws.Range("J62").Value = 5.49
$ws.Range("K62").Value = "17/09/2023 01:02"
$ws.Range("L62").Value = 7.5
$ws.Range("M62").Value = "26/09/2023 21:13"
$ws.Range("N62").Value = 3.92
$ws.Range("O62").Value = "17/09/2023 01:02"
$ws.Range("P62").Value = 4.4
$ws.Range("Q62").Value = "26/09/2023 21:13"
$ws.Range("R62").Value = 1.67
$ws.Range("S62").Value = "17/09/2023 01:02"
$ws.Range("T62").Value = 1.5
$ws.Range("U62").Value = "26/09/2023 21:13"
$ws.Range("V62").Value = "https://www.betexplorer.com/football/spain/laliga/mallorca-barcelona/noQWUFYc/"
